# The checklist used to have a single item + checkbox per row (columns A:B).
# Fix up the sheet so each row holds five "item / checkbox" pairs
# (A/B, C/D, E/F, G/H, I/J) and row 4 carries the real reading assignments.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old single-item labels in column A (rows 1-3 and 5); row 4 gets
# its new label below.
$ws.Range("A1").Value = ""
$ws.Range("A2").Value = ""
$ws.Range("A3").Value = ""
$ws.Range("A5").Value = ""

# New placeholder text (a single space) for the 4 extra "item" columns
# (C, E, G, I) on the rows that don't carry real text.
$ws.Range("C1").Value = " "
$ws.Range("C2").Value = " "
$ws.Range("C3").Value = " "
$ws.Range("C5").Value = " "

$ws.Range("E1").Value = " "
$ws.Range("E2").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("E5").Value = " "

$ws.Range("G1").Value = " "
$ws.Range("G2").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("G5").Value = " "

$ws.Range("I1").Value = " "
$ws.Range("I2").Value = " "
$ws.Range("I3").Value = " "
$ws.Range("I5").Value = " "

# Row 4 - the real reading checklist, spread across the 5 item columns.
$ws.Range("A4").Value = "Read Chapters 1-3"
$ws.Range("C4").Value = "Read Chapters 4-6 "
$ws.Range("E4").Value = " Read Chapters 7-8"
$ws.Range("G4").Value = " Read Chapters 9-10"
$ws.Range("I4").Value = " Read Chapters 11-13"

# The sheet's grid is 5 item/checkbox pairs wide (through column J) even
# though the last checkbox column never got a value - touch it so the
# sheet's used range extends out to J, matching the original template.
$ws.Range("J5").Font.Bold = $false
